# Updates odds values in rows 2, 7, 9, 10, 11 of the FlashScore weekly games sheet
# per the commit 'Atualizando o arquivo XLSX' (updating odds data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.67
$ws.Range("I2").Value = 4.33
$ws.Range("T2").Value = 13
$ws.Range("U2").Value = 12
$ws.Range("W2").Value = 17
$ws.Range("Y2").Value = 19
$ws.Range("AB2").Value = 13
$ws.Range("AF2").Value = 29

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 3.05
$ws.Range("L7").Value = 1.37
$ws.Range("M7").Value = 2.62
$ws.Range("N7").Value = 2.07
$ws.Range("P7").Value = 1.47
$ws.Range("Q7").Value = 2.32
$ws.Range("R7").Value = 1.82
$ws.Range("S7").Value = 1.78
$ws.Range("T7").Value = 6.7
$ws.Range("U7").Value = 10.5
$ws.Range("V7").Value = 9.25
$ws.Range("W7").Value = 23
$ws.Range("X7").Value = 21
$ws.Range("Y7").Value = 35
$ws.Range("Z7").Value = 7.8
$ws.Range("AA7").Value = 5.9
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 80
$ws.Range("AD7").Value = 800
$ws.Range("AF7").Value = 15
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 40
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 40

# Row 9
$ws.Range("G9").Value = 5.3
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 1.57
$ws.Range("M9").Value = 3.65
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 18
$ws.Range("U9").Value = 37
$ws.Range("V9").Value = 16.5
$ws.Range("W9").Value = 100
$ws.Range("X9").Value = 45
$ws.Range("Y9").Value = 40
$ws.Range("Z9").Value = 13
$ws.Range("AB9").Value = 14
$ws.Range("AC9").Value = 55
$ws.Range("AD9").Value = 350
$ws.Range("AE9").Value = 7.8
$ws.Range("AF9").Value = 8
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 11.75

# Row 10
$ws.Range("G10").Value = 1.38
$ws.Range("H10").Value = 4.75
$ws.Range("I10").Value = 7.5
$ws.Range("T10").Value = 9.5
$ws.Range("U10").Value = 8
$ws.Range("Z10").Value = 17
$ws.Range("AA10").Value = 9.5
$ws.Range("AE10").Value = 23
$ws.Range("AG10").Value = 21
$ws.Range("AI10").Value = 41

# Row 11
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 2.9
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 13
$ws.Range("L11").Value = 1.2
$ws.Range("M11").Value = 4.33
$ws.Range("N11").Value = 1.7
$ws.Range("O11").Value = 2.1
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.57
$ws.Range("S11").Value = 2.25
$ws.Range("T11").Value = 10
$ws.Range("U11").Value = 13
$ws.Range("Y11").Value = 23
$ws.Range("Z11").Value = 13
$ws.Range("AA11").Value = 6.5
$ws.Range("AB11").Value = 12
$ws.Range("AC11").Value = 41
$ws.Range("AD11").Value = 126
$ws.Range("AE11").Value = 12
$ws.Range("AJ11").Value = 26
